# Final Changes - 22 June 2023
#
# The "RecentlyViewedListView" sheet gains a new list entry ("PDC Sample")
# inserted just before the existing "Recently Viewed" / "Recently Viewed
# Contacts" rows, and becomes the active/selected sheet (with a new
# selection anchor), while the previously active "ExternalContactSections"
# sheet reverts to its default (unselected) view.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RecentlyViewedListView")

# Make this the active sheet (mirrors activeTab moving from
# ExternalContactSections to RecentlyViewedListView, and tabSelected
# moving accordingly).
$ws.Activate()

# Insert a new row above the old row 7 ("Recently Viewed") and populate it,
# pushing the existing rows 7-8 down to rows 8-9.
$ws.Rows("7:7").Insert()
$ws.Range("A7").Value = "PDC Sample"

# Update the selection to match the new active cell.
$ws.Range("D16").Select()
